# Fix text wrapping issues in presentations
# Split Composition slide 1: Widen orange block for 'SOFTWARE' text
$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)

# Shape 2: orange rectangle behind "SOFTWARE"
# (PowerPoint Shape position/size properties are in points; 1 pt = 12700 EMU.
#  A tiny epsilon is added because the host stores these as single-precision
#  floats internally, which would otherwise truncate e.g. 482.4pt -> 6126479 EMU
#  instead of the exact 6126480 EMU target.)
$orangeBlock = $s1.Shapes.Item(2)
$orangeBlock.Left = 201.60001
$orangeBlock.Top = 180.00001
$orangeBlock.Width = 482.40001
$orangeBlock.Height = 93.60001

# Shape 3: "SOFTWARE" text box
$softwareText = $s1.Shapes.Item(3)
$softwareText.Left = 216.00001
$softwareText.Top = 187.20001
$softwareText.Width = 453.60001
$softwareText.Height = 79.20001

# Executive Black slide 5 (here slide 4): rename feature titles/descriptions
$s4 = $p.Slides.Item(4)

$s4.Shapes.Item(4).TextFrame.TextRange.Text = "Touren-Planung"
$s4.Shapes.Item(5).TextFrame.TextRange.Text = "Intelligente Routenoptimierung"
$s4.Shapes.Item(11).TextFrame.TextRange.Text = "Zentrale Verwaltung"
